# Agile-Product-Backlog.xlsx : "Update up to sprint 8"
#
# Fills in task-tracking data for rows 25-30 (the last two task groups of
# Sprint 7 / "Implement security features" -> CRUD user functionality),
# adds the three SUM() rollup formulas that were missing for rows
# 23/27 (and restates the zero-value rollups for every still-empty
# group below), and renames a couple of backlog task names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Row 23 (group header "Implement security features") - rollup formula
# ---------------------------------------------------------------------
$ws.Range("J23").Formula = "=SUM(J24:J26)"

# ---------------------------------------------------------------------
# New task-name text is written in the same order the original author
# entered it so the shared-string table comes out in matching order:
#   C25 CreateUser, C26 ReadUsers, C30 UpdateUser, C29 Delete, C28 Read,
#   C27 CRUD (group header, entered last).
# ---------------------------------------------------------------------
$ws.Range("C25").Value = "Implement CreateUser method in NewUserAccount web form"
$ws.Range("C26").Value = "Implement ReadUsers method in User Index web form"

# Rows 28 & 30 reuse the borderless "s=4" style already present on C13.
$ws.Range("C13").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = "Implement UpdateUser method  in User account web form "

$ws.Range("C29").Value = "Implement Delete User method in "

$ws.Range("C13").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = "Implement ReadUser method  in User account web form "

$ws.Range("C27").Value = "Implement CRUD user functionality"

# ---------------------------------------------------------------------
# Row 25 - finish filling in the rest of the row
# ---------------------------------------------------------------------
$ws.Range("G25").Value = 43848
$ws.Range("H25").Value = "High"
$ws.Range("J25").Value = 40

# ---------------------------------------------------------------------
# Row 26 - finish filling in the rest of the row
# ---------------------------------------------------------------------
$ws.Range("B26").Value = 17
$ws.Range("F26").Value = 43849
$ws.Range("G26").Value = 43854
$ws.Range("H26").Value = "High"
$ws.Range("J26").Value = 40

# ---------------------------------------------------------------------
# Row 27 (group header) - rollup formula
# ---------------------------------------------------------------------
$ws.Range("J27").Formula = "=SUM(J28:J30)"

# ---------------------------------------------------------------------
# Row 28 - finish filling in the rest of the row
# ---------------------------------------------------------------------
$ws.Range("B28").Value = 18
$ws.Range("F28").Value = 43855
$ws.Range("G28").Value = 43857
$ws.Range("H28").Value = "High"
$ws.Range("J28").Value = 40

# ---------------------------------------------------------------------
# Row 29 - finish filling in the rest of the row
# ---------------------------------------------------------------------
$ws.Range("B29").Value = 19
$ws.Range("F29").Value = 43858
$ws.Range("G29").Value = 43861
$ws.Range("H29").Value = "High"
$ws.Range("J29").Value = 40

# ---------------------------------------------------------------------
# Row 30 - finish filling in the rest of the row
# ---------------------------------------------------------------------
$ws.Range("B30").Value = 20
$ws.Range("F30").Value = 43862
$ws.Range("G30").Value = 43866
$ws.Range("H30").Value = "High"
$ws.Range("J30").Value = 40

# ---------------------------------------------------------------------
# Rollup formulas for every other (still-empty) sprint group, rows
# 31-79. These all sum to 0 today since their child rows are blank.
# ---------------------------------------------------------------------
$ws.Range("J31").Formula = "=SUM(J32:J34)"
$ws.Range("J35").Formula = "=SUM(J36:J38)"
$ws.Range("J39").Formula = "=SUM(J40:J42)"
$ws.Range("J43").Formula = "=SUM(J44:J46)"
$ws.Range("J47").Formula = "=SUM(J48:J50)"
$ws.Range("J51").Formula = "=SUM(J52:J54)"
$ws.Range("J55").Formula = "=SUM(J56:J58)"
$ws.Range("J59").Formula = "=SUM(J60:J62)"
$ws.Range("J63").Formula = "=SUM(J64:J66)"
$ws.Range("J67").Formula = "=SUM(J68:J70)"
$ws.Range("J71").Formula = "=SUM(J72:J74)"
$ws.Range("J75").Formula = "=SUM(J76:J78)"
$ws.Range("J79").Formula = "=SUM(J80:J82)"

# ---------------------------------------------------------------------
# View state: scroll so the frozen pane starts around row 20 and the
# active selection lands on F30, matching where the author was working.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A20").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F30").Select()
